# 17Dec2021 Selenium DataDriven Part4
# Populate the "Result" column (S) of the TestData sheet with "PASS" for the
# rows that were previously left blank, matching the shared-string "PASS"
# (index 96) already used on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("S3").Value = "PASS"
$ws.Range("S5").Value = "PASS"
